# Insert a new data row at row 194, shifting the existing rows 194:277 down
# to 195:278, then populate the new row 194 with the added record.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows("194:194").Insert()

$ws.Range("A194").Value = 10
$ws.Range("B194").Value = 'Vega Modelo de Temuco'
$ws.Range("C194").Value = 'La Araucanía'
$ws.Range("D194").Value = 44609
$ws.Range("E194").Value = 9
$ws.Range("F194").Value = 100112009
$ws.Range("G194").Value = 'Acelga'
$ws.Range("H194").Value = 'Sin especificar'
$ws.Range("I194").Value = 'Primera'
$ws.Range("J194").Value = 50
$ws.Range("K194").Value = 8000
$ws.Range("L194").Value = 8000
$ws.Range("M194").Value = 8000
$ws.Range("N194").Value = '$/docena de atados (12 kilos)'
$ws.Range("O194").Value = 'Provincia de Cautín'
$ws.Range("P194").Value = 667
$ws.Range("Q194").Value = 12
$ws.Range("R194").Value = 'Hortaliza'
